$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = '@'
$ws.Range("D2").Value = '69.142.15'
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  +2.19%  '

$ws.Range("D3").NumberFormat = '@'
$ws.Range("D3").Value = '3.380.04'
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  +1.66%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = '@'
$ws.Range("D5").Value = '587.21'
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  +1.06%  '

$ws.Range("D6").NumberFormat = '@'
$ws.Range("D6").Value = '180.23'
$ws.Range("D6").Style = 'Normal'
$ws.Range("E6").Value = '  +2.78%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +0.90%  '

$ws.Range("D9").NumberFormat = '@'
$ws.Range("D9").Value = '0.197'
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +9.08%  '

$ws.Range("E10").Value = '  +1.51%  '

$ws.Range("D11").NumberFormat = '@'
$ws.Range("D11").Value = '48.67'
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  +4.83%  '

$ws.Range("E12").Value = '  +5.29%  '

$ws.Range("D13").NumberFormat = '@'
$ws.Range("D13").Value = '686.42'
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -2.58%  '

$ws.Range("E14").Value = '  +2.24%  '

$ws.Range("D15").NumberFormat = '@'
$ws.Range("D15").Value = '3.928.13'
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  +1.48%  '

$ws.Range("D16").NumberFormat = '@'
$ws.Range("D16").Value = '69.211.46'
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  +2.25%  '

$ws.Range("D17").NumberFormat = '@'
$ws.Range("D17").Value = '3.393.04'
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  +1.74%  '

$ws.Range("E18").Value = '  +1.82%  '

$ws.Range("D19").NumberFormat = '@'
$ws.Range("D19").Value = '17.71'
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  +2.09%  '

$ws.Range("E20").Value = '  +3.34%  '

$ws.Range("E21").Value = '  +0.81%  '

$ws.Range("E22").Value = '  +0.73%  '

$ws.Range("D23").NumberFormat = '@'
$ws.Range("D23").Value = '17.05'
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  +0.75%  '

$ws.Range("D24").NumberFormat = '@'
$ws.Range("D24").Value = '104.59'
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  +6.09%  '

$ws.Range("E25").Value = '  +1.53%  '

$ws.Range("E26").Value = '  +1.45%  '

$ws.Range("D27").NumberFormat = '@'
$ws.Range("D27").Value = '9.60'
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  +1.68%  '

$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '34.33'
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  +3.86%  '

$ws.Range("D29").NumberFormat = '@'
$ws.Range("D29").Value = '8.67'
$ws.Range("D29").Style = 'Normal'
$ws.Range("E29").Value = '  +1.78%  '

$ws.Range("D30").NumberFormat = '@'
$ws.Range("D30").Value = '6.97'
$ws.Range("D30").Style = 'Normal'
$ws.Range("E30").Value = '  -1.70%  '

$ws.Range("D31").NumberFormat = '@'
$ws.Range("D31").Value = '11.19'
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  +1.98%  '

$ws.Range("D32").NumberFormat = '@'
$ws.Range("D32").Value = '556.36'
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  -2.20%  '

$ws.Range("E33").Value = '  +10.18%  '

$ws.Range("D35").NumberFormat = '@'
$ws.Range("D35").Value = '57.95'
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  +0.79%  '

$ws.Range("E36").Value = '  +0.18%  '

$ws.Range("D37").NumberFormat = '@'
$ws.Range("D37").Value = '3.700.03'
$ws.Range("D37").Style = 'Normal'
$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("E38").Value = '  +7.55%  '

$ws.Range("D39").NumberFormat = '@'
$ws.Range("D39").Value = '34.84'
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  +2.21%  '

$ws.Range("D40").NumberFormat = '@'
$ws.Range("D40").Value = '3.23'
$ws.Range("D40").Style = 'Normal'
$ws.Range("E40").Value = '  +1.39%  '

$ws.Range("E41").Value = '  +4.63%  '

$ws.Range("E42").Value = '  +1.77%  '

$ws.Range("D43").NumberFormat = '@'
$ws.Range("D43").Value = '0.338'
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  +1.11%  '

$ws.Range("E44").Value = '  +3.21%  '

$ws.Range("E45").Value = '  -1.27%  '

$ws.Range("D46").NumberFormat = '@'
$ws.Range("D46").Value = '2.64'
$ws.Range("D46").Style = 'Normal'
$ws.Range("E46").Value = '  -0.87%  '

$ws.Range("E47").Value = '  +0.90%  '

$ws.Range("E48").Value = '  +5.00%  '

$ws.Range("E49").Value = '  -0.06%  '

$ws.Range("D50").NumberFormat = '@'
$ws.Range("D50").Value = '132.47'
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  +2.82%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = '@'
$ws.Range("D51").Value = '7.50'
$ws.Range("D51").Style = 'Normal'
$ws.Range("E51").Value = '  +1.25%  '

